$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.828.68"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "2.311.83"
$ws.Range("E3").Value = "  +1.57%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.ClearFormats()
$ws.Range("E4").Value = "  -0.15%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "301.76"
$c.ClearFormats()
$ws.Range("E5").Value = "  -1.23%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "96.47"
$c.ClearFormats()
$ws.Range("E6").Value = "  +0.07%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.509"
$c.ClearFormats()
$ws.Range("E7").Value = "  +0.47%  "
$ws.Range("E8").Value = "  -0.03%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.495"
$c.ClearFormats()
$ws.Range("E9").Value = "  -0.77%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "34.79"
$c.ClearFormats()
$ws.Range("E10").Value = "  -1.79%  "
$ws.Range("E11").Value = "  +5.72%  "
$ws.Range("E12").Value = "  -0.02%  "
$ws.Range("E13").Value = "  +0.12%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.81"
$c.ClearFormats()
$ws.Range("E14").Value = "  +1.26%  "
$ws.Range("D15").Value = "2.662.74"
$ws.Range("E15").Value = "  +1.09%  "
$ws.Range("D16").Value = "2.312.39"
$ws.Range("E16").Value = "  +1.15%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.788"
$c.ClearFormats()
$ws.Range("E17").Value = "  +1.49%  "
$ws.Range("D18").Value = "42.742.03"
$ws.Range("E18").Value = "  +0.33%  "
$ws.Range("E19").Value = "  -5.18%  "
$ws.Range("D20").Value = "0.0₃0894"
$ws.Range("E20").Value = "  -0.27%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "6.04"
$c.ClearFormats()
$ws.Range("E21").Value = "  +0.85%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "67.87"
$c.ClearFormats()
$ws.Range("E22").Value = "  +1.24%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "2.29"
$c.ClearFormats()
$ws.Range("E23").Value = "  +7.90%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "235.52"
$c.ClearFormats()
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("E26").Value = "  -1.74%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "24.47"
$c.ClearFormats()
$ws.Range("E27").Value = "  -2.56%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "2.36"
$c.ClearFormats()
$ws.Range("E28").Value = "  +15.19%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "165.41"
$c.ClearFormats()
$ws.Range("E29").Value = "  -0.37%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "9.09"
$c.ClearFormats()
$ws.Range("E30").Value = "  +0.66%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "32.32"
$c.ClearFormats()
$ws.Range("E31").Value = "  -1.91%  "
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("E33").Value = "  +0.62%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "17.67"
$c.ClearFormats()
$ws.Range("E34").Value = "  +0.78%  "
$ws.Range("E35").Value = "  -5.78%  "
$ws.Range("E36").Value = "  +2.24%  "
$ws.Range("E37").Value = "  -2.80%  "
$ws.Range("E38").Value = "  -0.74%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "1.76"
$c.ClearFormats()
$ws.Range("E39").Value = "  +1.27%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.72"
$c.ClearFormats()
$ws.Range("E40").Value = "  +0.94%  "
$ws.Range("E41").Value = "  -0.47%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "20.32"
$c.ClearFormats()
$ws.Range("E42").Value = "  +11.95%  "
$ws.Range("D43").Value = "1.976.88"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "10.46"
$c.ClearFormats()
$ws.Range("E44").Value = "  +5.24%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0280"
$c.ClearFormats()
$ws.Range("E45").Value = "  +0.46%  "
$ws.Range("E46").Value = "  -2.91%  "
$ws.Range("E47").Value = "  +0.80%  "
$ws.Range("D48").Value = "2.531.26"
$ws.Range("E48").Value = "  +1.11%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "53.50"
$c.ClearFormats()
$ws.Range("E49").Value = "  +0.05%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "2.76"
$c.ClearFormats()
$ws.Range("E50").Value = "  -2.76%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "71.63"
$c.ClearFormats()
$ws.Range("E51").Value = "  +0.48%  "
